$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (volume number + report week dates) ---
$ws.Range("A8").Characters(21,1).Text = "2"
$ws.Range("C9").Characters(48,8).Text = "1/12/2025"
$ws.Range("C9").Characters(27,10).Text = "1/6/2025"

# --- Anchor cells (format sources) kept stable elsewhere on sheet ---
# C33 s=13 t=s "0"(v=20); E33 s=13 t=s "***.*"(v=21); C39 s=14; K39 s=15

# --- Cells changing FROM text-placeholder TO numeric (need style+value) ---
$ws.Range("C39").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = 1
$ws.Range("K39").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = -100
$ws.Range("C39").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value = 1
$ws.Range("K39").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = -100
$ws.Range("C39").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("J14").Value = 1
$ws.Range("K39").Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("K14").Value = -100
$ws.Range("C39").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("C39").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("C39").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I15").Value = 1
$ws.Range("C39").Copy()
$ws.Range("J18").PasteSpecial(-4122)
$ws.Range("J18").Value = 4
$ws.Range("K39").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("K18").Value = 200
$ws.Range("K39").Copy()
$ws.Range("L22").PasteSpecial(-4122)
$ws.Range("L22").Value = 0
$ws.Range("C39").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 1
$ws.Range("C39").Copy()
$ws.Range("I23").PasteSpecial(-4122)
$ws.Range("I23").Value = 1
$ws.Range("K39").Copy()
$ws.Range("L23").PasteSpecial(-4122)
$ws.Range("L23").Value = 0
$ws.Range("C39").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("C39").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("C39").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("I27").Value = 1
$ws.Range("C39").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 2
$ws.Range("C39").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I28").Value = 2
$ws.Range("C39").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("K39").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("C39").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("K39").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("C39").Copy()
$ws.Range("J29").PasteSpecial(-4122)
$ws.Range("J29").Value = 1
$ws.Range("K39").Copy()
$ws.Range("K29").PasteSpecial(-4122)
$ws.Range("K29").Value = -100
$ws.Range("C39").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("K39").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("C39").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("K39").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100
$ws.Range("C39").Copy()
$ws.Range("J30").PasteSpecial(-4122)
$ws.Range("J30").Value = 1
$ws.Range("K39").Copy()
$ws.Range("K30").PasteSpecial(-4122)
$ws.Range("K30").Value = -100

# --- Cells changing FROM numeric TO text-placeholder (need style+value) ---
$ws.Range("C33").Copy()
$ws.Range("D20").PasteSpecial(-4104)
$ws.Range("E33").Copy()
$ws.Range("E20").PasteSpecial(-4104)
$ws.Range("C33").Copy()
$ws.Range("C22").PasteSpecial(-4104)
$ws.Range("C33").Copy()
$ws.Range("D22").PasteSpecial(-4104)
$ws.Range("E33").Copy()
$ws.Range("E22").PasteSpecial(-4104)
$ws.Range("C33").Copy()
$ws.Range("D28").PasteSpecial(-4104)
$ws.Range("E33").Copy()
$ws.Range("E28").PasteSpecial(-4104)
$ws.Range("C33").Copy()
$ws.Range("D31").PasteSpecial(-4104)
$ws.Range("E33").Copy()
$ws.Range("E31").PasteSpecial(-4104)

# --- Cells staying numeric, value-only change ---
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 9.090909090909
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 50
$ws.Range("L16").Value = 20
$ws.Range("M16").Value = -40
$ws.Range("N16").Value = -85
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -23.529411764705
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = -77.777777777777
$ws.Range("L17").Value = -81.818181818181
$ws.Range("N17").Value = -90
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 188.888888888889
$ws.Range("I18").Value = 12
$ws.Range("L18").Value = 500
$ws.Range("N18").Value = -45.454545454545
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 4.651162790697
$ws.Range("I19").Value = 18
$ws.Range("J19").Value = 19
$ws.Range("K19").Value = -5.263157894736
$ws.Range("L19").Value = 38.461538461538
$ws.Range("M19").Value = 5.882352941176
$ws.Range("N19").Value = -47.058823529411
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 3
$ws.Range("K20").Value = 200
$ws.Range("L20").Value = 200
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -86.95652173913
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -4.347826086956
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 17.045454545454
$ws.Range("I21").Value = 42
$ws.Range("J21").Value = 38
$ws.Range("K21").Value = 10.526315789473
$ws.Range("L21").Value = 31.25
$ws.Range("N21").Value = -69.78417266187
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 300
$ws.Range("M22").Value = -83.333333333333
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = -50
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 12.5
$ws.Range("F24").Value = 113
$ws.Range("G24").Value = 143
$ws.Range("H24").Value = -20.979020979021
$ws.Range("I24").Value = 49
$ws.Range("J24").Value = 51
$ws.Range("K24").Value = -3.92156862745
$ws.Range("L24").Value = -3.92156862745
$ws.Range("M24").Value = -3.92156862745
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -10.714285714285
$ws.Range("F25").Value = 84
$ws.Range("G25").Value = 117
$ws.Range("H25").Value = -28.205128205128
$ws.Range("I25").Value = 34
$ws.Range("J25").Value = 41
$ws.Range("K25").Value = -17.073170731707
$ws.Range("L25").Value = -22.727272727272
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 31
$ws.Range("H26").Value = 19.230769230769
$ws.Range("I26").Value = 16
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = 33.333333333333
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 33.333333333333
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -50
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("G31").Value = 1
$ws.Range("J42").Value = 338
$ws.Range("K42").Value = 69
$ws.Range("L42").Value = 20.714285714285
$ws.Range("M42").Value = -21.759259259259
$ws.Range("N42").Value = -11.979166666666
$ws.Range("J46").Value = 1713
$ws.Range("K46").Value = 1.722090261282
$ws.Range("L46").Value = -20.547309833024
$ws.Range("M46").Value = -61.660698299015
$ws.Range("N46").Value = -73.787299158378
